# ---------------------------------------------------------------------
# einfordern_schlussbemerkung.docx - "rsta templates ready for testing"
# ---------------------------------------------------------------------
$d = $word.ActiveDocument

# 1) Frame1 text-box: tiny re-layout / resize (wp:extent + a:ext + VML
#    fallback all shift by a fraction of a point). The Shape object's
#    Width/Height setters drive the drawing extents, so push them to the
#    new outer size.
$frame = $d.Shapes.Item("Frame1")
$frame.Width  = 2519045 / 12700.0
$frame.Height = 1798955 / 12700.0

# 2) "Oder Vertreter :" -> "Oder Vertreter:" (drop the space before the colon)
$d.Content.Find.Execute("Oder Vertreter :", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Oder Vertreter:", 2) | Out-Null

# 3) "{{VERTRETER_NAME_ADDRESS}} " -> " {{VERTRETER_NAME_ADDRESS}} " (leading space)
$d.Content.Find.Execute("{{VERTRETER_NAME_ADDRESS}} ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " {{VERTRETER_NAME_ADDRESS}} ", 2) | Out-Null

# 4) "Wiederholung Bet_1 gelöscht" paragraph: bump highlight cyan/yellow -> darkCyan
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Wiederholung Bet_1 gel*") {
        $p.Range.HighlightColorIndex = 10   # wdTeal -> w:highlight="darkCyan"
    }
}

# 5) Normal style: re-stamp the language so wd:val is (re)written first,
#    matching the canonical attribute order (no value change).
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.Font.LanguageID = "de-CH"

# 6) Register the next batch of (unused) list-label character styles,
#    continuing the ListLabel1..ListLabel80 sequence already in the doc.
$listLabelFonts = @(
    "Arial",       # 81
    "Courier New", # 82
    "Wingdings",   # 83
    "Symbol",      # 84
    "Courier New", # 85
    "Wingdings",   # 86
    "Symbol",      # 87
    "Courier New", # 88
    "Wingdings",   # 89
    "Arial",       # 90
    "Courier New", # 91
    "Wingdings",   # 92
    "Symbol",      # 93
    "Courier New", # 94
    "Wingdings",   # 95
    "Symbol",      # 96
    "Courier New", # 97
    "Wingdings"    # 98
)

for ($i = 0; $i -lt $listLabelFonts.Count; $i++) {
    $num = 81 + $i
    $newStyle = $d.Styles.Add("ListLabel " + $num, 2)
    $newStyle.Font.NameBi = $listLabelFonts[$i]
    $newStyle.QuickStyle = $true
}
